$d = $word.ActiveDocument

# 1) Title: "Unidad 1:" -> "Unidad 4:"
$d.Content.Find.Execute("Unidad 1:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Unidad 4:", 2) | Out-Null

# 2) Subtitle: " Algoritmos" -> " El entorno para la programación en java"
$d.Content.Find.Execute(" Algoritmos", $true, $false, $false, $false, $false,
                         $true, 1, $false, " El entorno para la programación en java", 2) | Out-Null

# 3) Merge date runs: "6" + " de mayo de 2021" -> single run "6 de mayo de 2021"
$d.Content.Find.Execute("6 de mayo de 2021", $true, $false, $false, $false, $false,
                         $true, 1, $false, "6 de mayo de 2021", 2) | Out-Null

# 4) Merge "Aprendizaje:" runs " " + "tipos de operadores en java, " -> " tipos de operadores en java, "
$d.Content.Find.Execute(" tipos de operadores en java, ", $true, $false, $false, $false, $false,
                         $true, 1, $false, " tipos de operadores en java, ", 2) | Out-Null

# 5) Merge "Reflexión: " runs " " + "En la clase de hoy..." -> " En la clase de hoy..."
$d.Content.Find.Execute(" En la clase de hoy vimos cuales eran los diferentes operadores que se pueden utilizar en java y es la parte que más me impactó debido a que si bien sé programar, nunca había visto lo comprimido que puede quedar un código. Aunque si bien queda comprimido, es un poquito más complejo de ver su funcionalidad. Por otra parte, vimos cómo funciona los ciclos (", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, `
                         " En la clase de hoy vimos cuales eran los diferentes operadores que se pueden utilizar en java y es la parte que más me impactó debido a que si bien sé programar, nunca había visto lo comprimido que puede quedar un código. Aunque si bien queda comprimido, es un poquito más complejo de ver su funcionalidad. Por otra parte, vimos cómo funciona los ciclos (", 2) | Out-Null
